# Applies the two edits described by the commit diff:
#
# 1. The opening paragraph ("This is a Microsoft word document.") gets two
#    trailing spaces appended to its existing run, then a new run in dark
#    red (C00000) containing "(This is a change – Version for branch
#    alternate)" is appended right after it.
#
# 2. The final "Bop-bop-bop-bop" at the very end of the lyrics paragraph
#    (the one that currently shares a run with the preceding line break)
#    gets isolated into its own run and wrapped in spell-check proofErr
#    markers, matching the shape of the other isolated "Bop-bop-bop-bop"
#    runs already present earlier in the same paragraph.

$d = $word.ActiveDocument

# --- Edit 1: opening sentence -------------------------------------------

$d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$firstPara = $d.Paragraphs(1).Range
$insertionPoint = $d.Range($firstPara.End - 1, $firstPara.End - 1)
$insertionPoint.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
$insertionPoint.Font.Color = 192   # wdColor is BGR -> 0x0000C0 == RGB C00000

# --- Edit 2: isolate the trailing "Bop-bop-bop-bop" ----------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$targetWord = "Bop-bop-bop-bop"
$wordEnd = $lastPara.End - 1
$wordStart = $wordEnd - $targetWord.Length
$target = $d.Range($wordStart, $wordEnd)

# Re-asserting the (already-effective) East Asian font name forces the
# run to split off from the preceding line-break run while keeping every
# other character property identical; the engine surrounds the newly
# isolated run with spell-check proofErr markers automatically.
$target.Font.NameFarEast = "Times New Roman"
